$d = $word.ActiveDocument

# --- Step 1: capture the formatted text of the existing
#     "{% if loop.index == 1 %} and {% endif %}" block so the new copy
#     keeps identical run/formatting structure.
$src = $d.Content
$src.Find.Execute("{% if loop.index == 1 %} and {% endif %}", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$ft = $src.FormattedText

# --- Step 2: paste a copy of that block right after
#     "{% for charge in amended_charges_list %}" (i.e. right before
#     "{{ charge[0] }}").
$dest = $d.Content
$dest.Find.Execute("amended_charges_list %}", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$dest.Collapse(0)  # wdCollapseEnd
$insertStart = $dest.Start
$dest.FormattedText = $ft

# The newly inserted copy currently reads "{% if loop.index == 1 %} and {% endif %}".
# Change its "1" to "2" -> "{% if loop.index == 2 %} and {% endif %}"
$newBlock = $d.Range($insertStart, $insertStart + $ft.Text.Length)
$newBlock.Find.Execute("loop.index == 1", $false, $false, $false, $false, $false, $true, 1, $false, "loop.index == 2", 2) | Out-Null

# --- Step 3: put the _GoBack bookmark at this new edit location
#     (matching Word's habit of tracking the most recent edit point),
#     removing it from wherever it previously sat.
$old = $d.Bookmarks.Item("_GoBack")
$old.Delete()
$d.Bookmarks.Add("_GoBack", $d.Range($insertStart, $insertStart))

# --- Step 4: remove the old "{% if loop.index == 1 %} and {% endif %}"
#     that used to sit between "{{ charge[1] }}" and "{% endfor %}".
#     (Find again from the top so the match after our insertion is the
#     original one, further along in the document.)
$old2 = $d.Content
$old2.Find.Execute("{% if loop.index == 1 %} and {% endif %}", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$old2.Text = ""

# --- Step 5: collapse the trailing double endif down to one.
#     "...{% endif %} {% endif %}" -> "...{% endif %} {% endif %}"
#     Remove the first of the pair (keep the leading space before the
#     second).
$tail = $d.Content
$tail.Find.Execute("{% endif %} {% endif %}", $false, $false, $false, $false, $false, $true, 1, $false, " {% endif %}", 2) | Out-Null
